{"js": "const body = context.document.body;\n\n// --- Change 1 -------------------------------------------------------------\n// \"...von 0 bis 50 reichen.\" becomes\n// \"...von 0 bis 1000 (space_size, 1000 vielleicht nicht endg\u00fcltig) reichen.\"\nconst sizeResults = body.search(\"50\", { matchCase: true, matchWholeWord: true });\nsizeResults.load(\"items\");\nawait context.sync();\n\nif (sizeResults.items.length === 0) {\n  throw new Error('Could not find \"50\" to expand into the new coordinate-range text.');\n}\n\nconst sizeRange = sizeResults.items[0];\nconst insertedRange = sizeRange.insertText(\n  \"1000 (space_size, 1000 vielleicht nicht endg\u00fcltig)\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Change 2 -------------------------------------------------------------\n// The \"_GoBack\" bookmark moves from the end of the \"Die Funktion gibt...\"\n// paragraph to right after the text that was just inserted above.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst bookmarkSpot = insertedRange.getRange(Word.RangeLocation.end);\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Change 3 -------------------------------------------------------------\n// A new paragraph about the \"Randknoten\" is appended after the paragraph\n// ending in \"...die zweite y-Werte.\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst afterParagraph = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Die Funktion gibt einen zweidimensionalen Array zur\u00fcck\") !== -1\n);\nif (!afterParagraph) {\n  throw new Error(\"Could not find the paragraph to insert the new Randknoten paragraph after.\");\n}\n\nafterParagraph.insertParagraph(\n  \"Die Randknoten werden einzeln generiert, sie nehmen mindestens einen Extremwert an. (Anmerkung: Da die Randknoten einzeln in einer Funktion generiert werden, wird es nachher einfacher sein, die Randknoten mit einer h\u00f6heren spawn-rate zu versehen.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 ---------------------------------------------------------\n# \"...von 0 bis 50 reichen.\" becomes\n# \"...von 0 bis 1000 (space_size, 1000 vielleicht nicht endg\u00fcltig) reichen.\"\n$sizeRange = $d.Content\n$find = $sizeRange.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$found = $find.Execute(\n    \"50\",\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"1000 (space_size, 1000 vielleicht nicht endg\u00fcltig)\",\n    2\n)\nif (-not $found) {\n    throw 'Could not find \"50\" to expand into the new coordinate-range text.'\n}\n\n# --- Change 2 -----------------------------------------------------------\n# The \"_GoBack\" bookmark moves from the end of the \"Die Funktion gibt...\"\n# paragraph to right after the text that was just inserted above.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n$sizeRange.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $sizeRange)\n\n# --- Change 3 -------------------------------------------------------------\n# A new paragraph about the \"Randknoten\" is appended after the paragraph\n# ending in \"...die zweite y-Werte.\"\n$paragraphs = $d.Paragraphs\n$afterParagraph = $paragraphs.Item($paragraphs.Count)\n$afterParagraph.Range.InsertParagraphAfter()\n$paragraphs = $d.Paragraphs\n$newParagraph = $paragraphs.Item($paragraphs.Count)\n$newParagraph.Range.Text = \"Die Randknoten werden einzeln generiert, sie nehmen mindestens einen Extremwert an. (Anmerkung: Da die Randknoten einzeln in einer Funktion generiert werden, wird es nachher einfacher sein, die Randknoten mit einer h\u00f6heren spawn-rate zu versehen.\"\n"}
